# Update "想去人数" (number of interested attendees) values on the
# "展览" and "全部类型" sheets, as produced by the site regeneration.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F5").Value = 12956
    $ws.Range("F12").Value = 13729
    $ws.Range("F13").Value = 14231
    $ws.Range("F25").Value = 937
    $ws.Range("F26").Value = 5313
}
